# Auto Mode List.xlsx - Queen City Regional update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# The list shrank from 14 data rows to 9: drop the old rows 10-14 entirely
$ws.Range("A10:A14").EntireRow.Delete()

# Row 1 (title banner) and Row 2 (column headers) keep their existing formatting;
# only the row-2 values need to flip to the new shared-string order, which happens
# automatically because the text itself does not change ("MODE #" / "Defense").

# --- Build the new green-on-black "data row" look for column A (centered) ---
# Start from an existing green-on-black cell (B4) so the font/fill are reused,
# then add centered alignment -> becomes the new style used for A3:A9.
$ws.Range("B4").Copy()
$ws.Range("A3:A9").PasteSpecial(-4122)
$ws.Range("A3:A9").HorizontalAlignment = -4108

# Column B green-on-black, left/general aligned, reusing the existing style from B4
$ws.Range("B4").Copy()
$ws.Range("B3:B9").PasteSpecial(-4122)

# Now fill in the new values
$ws.Range("A3").Value = 8
$ws.Range("B3").Value = "Corner Shot"

$ws.Range("A4").Value = 100
$ws.Range("B4").Value = "Low Bar One Ball (w Gyro)"

$ws.Range("A5").Value = 101
$ws.Range("B5").Value = "Portcullis One Ball (w Gyro)"

$ws.Range("A6").Value = 102
$ws.Range("B6").Value = "Cheval One Ball (w Gyro)"

$ws.Range("A7").Value = 103
$ws.Range("B7").Value = "Rough Terrain One Ball (w Gyro)"

$ws.Range("A8").Value = 200
$ws.Range("B8").Value = "Low Bar Two Ball  w Spybot (w Gyro)"

$ws.Range("A9").Value = "default"
$ws.Range("B9").Value = "Corner Shot"

# Row 8 (the call-out "200" two-ball auto) gets yellow-on-black instead of green
$ws.Range("A8").Font.Color = 65535
$ws.Range("B8").Font.Color = 65535

$ws.Range("A1:B9").Select()
